$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps: Hedera <-> InternetComputer(DFINITY); Decentraland <-> EnergySwap ---
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.414"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06614"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.96"
$ws.Range("E45").Value = "  +2.09%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6341"
$ws.Range("E46").Value = "  +2.45%  "

# --- Price / Volume(1h) updates ---
$ws.Range("D2").Value = "30.208.73"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.070.48"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.60"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5180"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4319"
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08908"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.72"
$ws.Range("E10").Value = "  +6.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.151"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.16"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "2.071.22"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.641"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.640"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.74"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001119"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06600"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.70"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "30.260.61"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.20"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").Value = "2.315.85"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.14"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.520"
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.41"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.186"
$ws.Range("E31").Value = "  +3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.620"
$ws.Range("E33").Value = "  +19.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.075"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.819"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02563"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.639"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.52"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2234"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6811"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.238"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9975"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.194"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.594"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.190"
$ws.Range("E50").Value = "  +7.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.02"
$ws.Range("E51").Value = "  -0.54%  "
